$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.650499999999999
$ws.Range("B6").Value = 5.628600000000002
$ws.Range("B7").Value = 5.023399999999999
$ws.Range("C7").Value = -14.0014
$ws.Range("C12").Value = -11.2275
$ws.Range("D13").Value = -8.302199999999997
$ws.Range("D14").Value = -7.766500000000001
$ws.Range("C15").Value = -14.85819999999999
$ws.Range("B16").Value = 7.514499999999999
$ws.Range("D16").Value = -8.688400000000007
$ws.Range("D19").Value = -8.304299999999991
$ws.Range("B20").Value = 9.402099999999997
$ws.Range("C20").Value = -11.87190000000001
$ws.Range("C21").Value = -12.00530000000002
$ws.Range("C22").Value = -13.29669999999999
$ws.Range("D22").Value = -7.918299999999999
$ws.Range("C23").Value = -12.11790000000001
$ws.Range("B28").Value = 5.9671
$ws.Range("B29").Value = 4.694999999999999
$ws.Range("C29").Value = -10.2691
$ws.Range("B32").Value = 7.096899999999994
$ws.Range("C34").Value = -11.52400000000002
$ws.Range("D36").Value = -8.355699999999993
$ws.Range("B40").Value = 9.427899999999992
$ws.Range("C42").Value = -12.6742
$ws.Range("C43").Value = -13.12009999999999
$ws.Range("C44").Value = -13.21439999999999
$ws.Range("C45").Value = -13.98609999999999
$ws.Range("B46").Value = 6.108099999999999
$ws.Range("C46").Value = -13.2412
$ws.Range("D46").Value = -8.180100000000003
$ws.Range("C50").Value = -14.17199999999999
$ws.Range("D50").Value = -7.9616
$ws.Range("B51").Value = 5.854300000000003
$ws.Range("C51").Value = -12.27250000000001
$ws.Range("B52").Value = 4.978900000000002
$ws.Range("B57").Value = 5.093699999999997
$ws.Range("B59").Value = 4.832199999999999
$ws.Range("B62").Value = 5.620699999999999
$ws.Range("B66").Value = 5.495100000000001
$ws.Range("C66").Value = -11.43980000000001
$ws.Range("C67").Value = -11.3303
$ws.Range("B73").Value = 8.372400000000003
$ws.Range("B74").Value = 9.23009999999999
$ws.Range("C79").Value = -11.3493
$ws.Range("C84").Value = -13.42139999999999
$ws.Range("B92").Value = 5.582499999999993
$ws.Range("C92").Value = -11.1995
$ws.Range("D95").Value = -8.174199999999997
$ws.Range("C97").Value = -11.6248
$ws.Range("D97").Value = -8.512699999999997
$ws.Range("B100").Value = 5.825999999999997
